# "retrieve sheet dan read by sheet"
# Adds a new "GeneralMasterLoan" worksheet (after the existing
# "GeneralMasterEmployee" sheet), fills it with loan/nasabah data, formats
# the "Jumlah Hutang" column with a thousands-separator number format, and
# makes the new sheet the active one (matching tabSelected / activeTab in
# the target workbook).

$wb = $excel.ActiveWorkbook

# --- add the new worksheet after the last existing sheet -------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "GeneralMasterLoan"

# --- column widths (best-effort match of the authored widths) --------------
$ws.Columns.Item(2).ColumnWidth = 16.33
$ws.Columns.Item(3).ColumnWidth = 10
$ws.Columns.Item(4).ColumnWidth = 17.67

# --- data (id, nama nasabah, kota tinggal, jumlah hutang, status, dpd) ----
$ids      = @(1, 2, 3, 4, 5, 6, 7, 8, 9, 10)
$names    = @("Eko Patrio", "Sumanto", "Alex", "Baharuddin", "Suprapto", "Raharjo", "Bambang", "Suminto", "Ponaryo", "Reza")
$cities   = @("Semarang", "Denpasar", "Bandung", "Solo", "Jakarta", "Semarang", "Solo", "Denpasar", "Bandung", "Surabaya")
$amounts  = @(5000000, 2500000, 3450000, 10000000, 4300000, 2100000, 1450000, 1200000, 450000, 500000)
$statuses = @("Nunggak", "Nunggak", "Nunggak", "Nunggak", "Nunggak", "Nunggak", "Nunggak", "Nunggak", "Nunggak", "Nunggak")
$dpds     = @(5, 4, 10, 20, 13, 20, 3, 4, 2, 1)

# Fill in the same order the sheet was originally authored in (headers for
# A:E first, then each data column top-to-bottom, then the DPD column last)
# so the resulting shared-string table indices line up.
$ws.Range("A2").Value = "Id"
$ws.Range("B2").Value = "Nama Nasabah"
$ws.Range("C2").Value = "Kota Tinggal"
$ws.Range("D2").Value = "Jumlah Hutang"
$ws.Range("E2").Value = "Status"

for ($i = 0; $i -lt $ids.Count; $i++) { $ws.Range("A$(3+$i)").Value = $ids[$i] }
for ($i = 0; $i -lt $names.Count; $i++) { $ws.Range("B$(3+$i)").Value = $names[$i] }
for ($i = 0; $i -lt $cities.Count; $i++) { $ws.Range("C$(3+$i)").Value = $cities[$i] }
for ($i = 0; $i -lt $amounts.Count; $i++) {
    $cell = $ws.Range("D$(3+$i)")
    $cell.Value = $amounts[$i]
    $cell.NumberFormat = "#,##0"
}
for ($i = 0; $i -lt $statuses.Count; $i++) { $ws.Range("E$(3+$i)").Value = $statuses[$i] }

$ws.Range("F2").Value = "DPD"
for ($i = 0; $i -lt $dpds.Count; $i++) { $ws.Range("F$(3+$i)").Value = $dpds[$i] }

# --- make the new sheet the active / selected tab -------------------------
$ws.Activate() | Out-Null
$ws.Range("F13").Select() | Out-Null
